$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.441.68"
$ws.Range("E2").Value = "  -3.47%  "

# Row 3
$ws.Range("D3").Value = "3.705.21"
$ws.Range("E3").Value = "  -5.81%  "

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "597.07"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "166.65"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.47%  "

# Row 7
$ws.Range("D7").Value = "3.702.72"
$ws.Range("E7").Value = "  -5.73%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.532"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.163"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.67%  "

# Row 11
$ws.Range("E11").Value = "  -3.76%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.463"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -4.68%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "37.83"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -5.57%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000242"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -5.46%  "

# Row 15
$ws.Range("D15").Value = "4.320.50"
$ws.Range("E15").Value = "  -5.78%  "

# Row 16
$ws.Range("D16").Value = "3.702.50"
$ws.Range("E16").Value = "  -5.93%  "

# Row 17
$ws.Range("D17").Value = "67.446.52"
$ws.Range("E17").Value = "  -3.58%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "7.26"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -3.23%  "

# Row 19
$ws.Range("E19").Value = "  +5.22%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.66%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "487.08"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -4.11%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.27"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -5.77%  "

# Row 23
$ws.Range("E23").Value = "  -4.00%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "85.27"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.16%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -6.33%  "

# Row 26
$ws.Range("E26").Value = "  -1.81%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "12.19"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.49%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.08"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.11%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.44%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.35"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -8.58%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.71"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.23%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "31.32"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -7.97%  "

# Row 34
$ws.Range("D34").Value = "3.844.70"
$ws.Range("E34").Value = "  -6.01%  "

# Row 35
$ws.Range("E35").Value = "  -5.04%  "

# Row 36
$ws.Range("D36").Value = "3.644.33"
$ws.Range("E36").Value = "  -5.57%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.85%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.84"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -5.35%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.131"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -7.93%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.322"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.87%  "

# Row 42
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "48.65"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "424.97"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -10.46%  "

# Row 44
$ws.Range("E44").Value = "  -5.50%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.72%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "8.45"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.76%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "40.43"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -5.42%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "140.85"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0351"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.96%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.738.89"
$ws.Range("E51").Value = "  -7.48%  "
